# Apply updated cryptocurrency price/volume data to sheet1 (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.586.51"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.995.21"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.81"
$ws.Range("E5").Value = "  -9.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.598"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.89"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.76"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.20"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "2.287.28"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("E15").Value = "  -2.62%  "
$ws.Range("E16").Value = "  -6.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.09"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("D18").Value = "1.995.86"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "36.493.60"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.90"
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").Value = "0.0₃0807"
$ws.Range("E21").Value = "  -3.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.28"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.91"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  -9.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.20"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.68"
$ws.Range("E28").Value = "  -3.90%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.87"
$ws.Range("E30").Value = "  -4.58%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -6.48%  "
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.27"
$ws.Range("E35").Value = "  -7.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -2.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.60"
$ws.Range("E40").Value = "  +6.08%  "
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0948"
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("D43").Value = "1.452.31"
$ws.Range("E43").Value = "  +3.31%  "
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -8.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.28"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.23"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.995"
$ws.Range("E48").Value = "  -3.40%  "
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.71"
$ws.Range("E51").Value = "  +7.13%  "
